# Append a new row (row 4) to Sheet1, mirroring the shape of the existing
# rows (2 and 3) with the new trip's data and timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4 is blank text (like A2/A3). Using a bare apostrophe forces an empty
# text cell instead of clearing it outright; resetting the style back to
# Normal drops the quote-prefix formatting flag so it matches A2/A3.
$ws.Range("A4").Value = "'"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = "احمد"

# C4 ("233") looks numeric, so force it to stay text the same way.
$ws.Range("C4").Value = "'233"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = "الصمود"
$ws.Range("E4").Value = "الرحلة 2"
$ws.Range("F4").Value = "C2"
$ws.Range("G4").Value = "IDRF"
$ws.Range("H4").Value = "٠٥‏/٠٥‏/٢٠٢٥ ٠٢:٠١:٢٠ م"
